# Add the new "optimization_parameters" worksheet as the last sheet and
# make it the active sheet (matches the workbook.xml bookViews/sheets diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "optimization_parameters"

# Row 1 - headers
$ws.Range("A1").Value = "optimization_parameter"
$ws.Range("B1").Value = "value"

# Row 2 - alpha
$ws.Range("A2").Value = "alpha"
$ws.Range("B2").Value = 0.002
$ws.Range("B2").NumberFormat = "0.00E+00"

# Row 3 - kk_max
$ws.Range("A3").Value = "kk_max"
$ws.Range("B3").Value = 1

# Row 4 - MaxIter
$ws.Range("A4").Value = "MaxIter"
$ws.Range("B4").Value = 100000000
$ws.Range("B4").NumberFormat = "0.00E+00"

# Row 5 - TolFun
$ws.Range("A5").Value = "TolFun"
$ws.Range("B5").Value = 0.000001
$ws.Range("B5").NumberFormat = "0.00E+00"

# Row 6 - MaxFunEval
$ws.Range("A6").Value = "MaxFunEval"
$ws.Range("B6").Value = 100000000
$ws.Range("B6").NumberFormat = "0.00E+00"

# Row 7 - TolX
$ws.Range("A7").Value = "TolX"
$ws.Range("B7").Value = 0.000001
$ws.Range("B7").NumberFormat = "0.00E+00"

# Row 8 - production_function
$ws.Range("A8").Value = "production_function"
$ws.Range("B8").Value = "Sigmoid"

# Row 9 - L_curve
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Row 10 - estimate_params
$ws.Range("A10").Value = "estimate_params"
$ws.Range("B10").Value = 1

# Row 11 - make_graphs
$ws.Range("A11").Value = "make_graphs"
$ws.Range("B11").Value = 1

# Row 12 - fix_P
$ws.Range("A12").Value = "fix_P"
$ws.Range("B12").Value = 0

# Row 13 - fix_b
$ws.Range("A13").Value = "fix_b"
$ws.Range("B13").Value = 0

# Row 14 - expression_timepoints
$ws.Range("A14").Value = "expression_timepoints"
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 60

# Row 15 - Strain
$ws.Range("A15").Value = "Strain"
$ws.Range("B15").Value = "wt"
$ws.Range("C15").Value = "dgln3"

# Row 16 - simulation_timepoints
$ws.Range("A16").Value = "simulation_timepoints"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 15

# Row 17 - species
$ws.Range("A17").Value = "species"
$ws.Range("B17").Value = "Saccharomyces cerevisiae"

# Row 18 - taxon_id
$ws.Range("A18").Value = "taxon_id"
$ws.Range("B18").Value = 559292

# Column width to roughly match the target sheet layout (first column is
# wider to fit the parameter names).
$ws.Columns.Item(1).ColumnWidth = 21.5

# Select A1 as the active cell on the new sheet, matching the target's
# selection anchor (closest match to the authored sheet's activeCell).
$ws.Range("A1").Select()
